$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-10-06 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-07 Tuesday", 2) | Out-Null

# Update each answer cell in the practice table (20 rows x 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "83-68=15"
$t.Cell(1, 2).Range.Text = "33+64=97"
$t.Cell(1, 3).Range.Text = "0+73=73"
$t.Cell(1, 4).Range.Text = "81-60=21"
$t.Cell(1, 5).Range.Text = "8+82=90"

$t.Cell(2, 1).Range.Text = "36+55=91"
$t.Cell(2, 2).Range.Text = "44+16=60"
$t.Cell(2, 3).Range.Text = "29+23=52"
$t.Cell(2, 4).Range.Text = "72-13=59"
$t.Cell(2, 5).Range.Text = "71-48=23"

$t.Cell(3, 1).Range.Text = "35+43=78"
$t.Cell(3, 2).Range.Text = "77-26=51"
$t.Cell(3, 3).Range.Text = "3+75=78"
$t.Cell(3, 4).Range.Text = "22+16=38"
$t.Cell(3, 5).Range.Text = "88-77=11"

$t.Cell(4, 1).Range.Text = "42+38=80"
$t.Cell(4, 2).Range.Text = "7+31=38"
$t.Cell(4, 3).Range.Text = "48-17=31"
$t.Cell(4, 4).Range.Text = "72-4=68"
$t.Cell(4, 5).Range.Text = "18+70=88"

$t.Cell(5, 1).Range.Text = "64-44=20"
$t.Cell(5, 2).Range.Text = "77+15=92"
$t.Cell(5, 3).Range.Text = "30-25=5"
$t.Cell(5, 4).Range.Text = "1+42=43"
$t.Cell(5, 5).Range.Text = "48-9=39"

$t.Cell(6, 1).Range.Text = "39-5=34"
$t.Cell(6, 2).Range.Text = "71+17=88"
$t.Cell(6, 3).Range.Text = "33+39=72"
$t.Cell(6, 4).Range.Text = "14-14=0"
$t.Cell(6, 5).Range.Text = "90-39=51"

$t.Cell(7, 1).Range.Text = "30+35=65"
$t.Cell(7, 2).Range.Text = "63-20=43"
$t.Cell(7, 3).Range.Text = "13+6=19"
$t.Cell(7, 4).Range.Text = "72-16=56"
$t.Cell(7, 5).Range.Text = "77-5=72"

$t.Cell(8, 1).Range.Text = "75-25=50"
$t.Cell(8, 2).Range.Text = "60-30=30"
$t.Cell(8, 3).Range.Text = "80-75=5"
$t.Cell(8, 4).Range.Text = "5+26=31"
$t.Cell(8, 5).Range.Text = "82-29=53"

$t.Cell(9, 1).Range.Text = "83-2=81"
$t.Cell(9, 2).Range.Text = "61-41=20"
$t.Cell(9, 3).Range.Text = "18+51=69"
$t.Cell(9, 4).Range.Text = "89-5=84"
$t.Cell(9, 5).Range.Text = "82-65=17"

$t.Cell(10, 1).Range.Text = "0+49=49"
$t.Cell(10, 2).Range.Text = "52-38=14"
$t.Cell(10, 3).Range.Text = "42-16=26"
$t.Cell(10, 4).Range.Text = "46-40=6"
$t.Cell(10, 5).Range.Text = "91-26=65"

$t.Cell(11, 1).Range.Text = "64+13=77"
$t.Cell(11, 2).Range.Text = "48-46=2"
$t.Cell(11, 3).Range.Text = "59-57=2"
$t.Cell(11, 4).Range.Text = "75-61=14"
$t.Cell(11, 5).Range.Text = "45-40=5"

$t.Cell(12, 1).Range.Text = "3+1=4"
$t.Cell(12, 2).Range.Text = "51-49=2"
$t.Cell(12, 3).Range.Text = "96-87=9"
$t.Cell(12, 4).Range.Text = "77-46=31"
$t.Cell(12, 5).Range.Text = "22-17=5"

$t.Cell(13, 1).Range.Text = "88-78=10"
$t.Cell(13, 2).Range.Text = "34+52=86"
$t.Cell(13, 3).Range.Text = "9+26=35"
$t.Cell(13, 4).Range.Text = "93-25=68"
$t.Cell(13, 5).Range.Text = "97-95=2"

$t.Cell(14, 1).Range.Text = "13+24=37"
$t.Cell(14, 2).Range.Text = "4+3=7"
$t.Cell(14, 3).Range.Text = "57+0=57"
$t.Cell(14, 4).Range.Text = "26+34=60"
$t.Cell(14, 5).Range.Text = "72+2=74"

$t.Cell(15, 1).Range.Text = "44+7=51"
$t.Cell(15, 2).Range.Text = "57+38=95"
$t.Cell(15, 3).Range.Text = "78-18=60"
$t.Cell(15, 4).Range.Text = "57-54=3"
$t.Cell(15, 5).Range.Text = "15-13=2"

$t.Cell(16, 1).Range.Text = "98-35=63"
$t.Cell(16, 2).Range.Text = "0+8=8"
$t.Cell(16, 3).Range.Text = "85-48=37"
$t.Cell(16, 4).Range.Text = "64+15=79"
$t.Cell(16, 5).Range.Text = "6+53=59"

$t.Cell(17, 1).Range.Text = "94-92=2"
$t.Cell(17, 2).Range.Text = "80-27=53"
$t.Cell(17, 3).Range.Text = "23+65=88"
$t.Cell(17, 4).Range.Text = "97-71=26"
$t.Cell(17, 5).Range.Text = "36+53=89"

$t.Cell(18, 1).Range.Text = "53-17=36"
$t.Cell(18, 2).Range.Text = "12-1=11"
$t.Cell(18, 3).Range.Text = "10+40=50"
$t.Cell(18, 4).Range.Text = "32-18=14"
$t.Cell(18, 5).Range.Text = "5+76=81"

$t.Cell(19, 1).Range.Text = "89-11=78"
$t.Cell(19, 2).Range.Text = "68-5=63"
$t.Cell(19, 3).Range.Text = "21-1=20"
$t.Cell(19, 4).Range.Text = "2+50=52"
$t.Cell(19, 5).Range.Text = "54+43=97"

$t.Cell(20, 1).Range.Text = "24+50=74"
$t.Cell(20, 2).Range.Text = "27+39=66"
$t.Cell(20, 3).Range.Text = "29-29=0"
$t.Cell(20, 4).Range.Text = "44-43=1"
$t.Cell(20, 5).Range.Text = "25+10=35"
